$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.07206204654897827
$ws.Range("C2").Value = -7.268604489308395
$ws.Range("D2").Value = 1.080217527560668
$ws.Range("B3").Value = 0.06360257439174907
$ws.Range("C3").Value = -5.821495374793575
$ws.Range("D3").Value = 0.6526024569378328
$ws.Range("B4").Value = 0.03185079623909565
$ws.Range("C4").Value = -4.636029275706818
$ws.Range("D4").Value = 0.4532677099663868
$ws.Range("B5").Value = 0.06239001471151637
$ws.Range("C5").Value = -3.547987541353218
$ws.Range("D5").Value = 0.4630013904310029
$ws.Range("B6").Value = 0.03180531377609405
$ws.Range("C6").Value = -2.920083701141636
$ws.Range("D6").Value = 0.3874397482008062
$ws.Range("B7").Value = 0.03262772940045123
$ws.Range("C7").Value = -2.243220563562845
$ws.Range("D7").Value = 0.2442271716550833
$ws.Range("B8").Value = 0.03524416060046009
$ws.Range("C8").Value = -1.619903500917138
$ws.Range("D8").Value = 0.1875941329790004
$ws.Range("B9").Value = 0.02578820514380926
$ws.Range("C9").Value = -1.080191826876142
$ws.Range("D9").Value = 0.1208157567186046
$ws.Range("B10").Value = 0.03860807742207144
$ws.Range("C10").Value = -0.5657305909400557
$ws.Range("D10").Value = 0.06691744012641497
$ws.Range("B11").Value = 0.03399631035228703
$ws.Range("C11").Value = -0.3026927446299457
$ws.Range("D11").Value = 0.03771603205386091
$ws.Range("B12").Value = 0.04405952956379752
$ws.Range("C12").Value = 0.08833350720392884
$ws.Range("D12").Value = 0.009784809329549084
$ws.Range("B13").Value = 0.03946391287327365
$ws.Range("C13").Value = 0.6336856394345912
$ws.Range("D13").Value = 0.07209904453772546
$ws.Range("B14").Value = 0.06730611017008925
$ws.Range("C14").Value = 1.14481553519387
$ws.Range("D14").Value = 0.08897525267211621
$ws.Range("B15").Value = 0.03567449926840217
$ws.Range("C15").Value = 1.398751566532662
$ws.Range("D15").Value = 0.1106277072034678
$ws.Range("B16").Value = 0.0353629329142244
$ws.Range("C16").Value = 2.050373643860772
$ws.Range("D16").Value = 0.1253059262376374
$ws.Range("B17").Value = 0.06107668597927941
$ws.Range("C17").Value = 2.789130981393789
$ws.Range("D17").Value = 0.3781280492398744
$ws.Range("B18").Value = 0.07005319487477937
$ws.Range("C18").Value = 3.74785824886251
$ws.Range("D18").Value = 0.4706154026733314
$ws.Range("B19").Value = 0.05932635130636818
$ws.Range("C19").Value = 4.463304770990232
$ws.Range("D19").Value = 0.3166609705400714
$ws.Range("B20").Value = 0.04823075933267564
$ws.Range("C20").Value = 5.843529599134729
$ws.Range("D20").Value = 0.736505663861722
$ws.Range("B21").Value = 0.03485852443155366
$ws.Range("C21").Value = 7.377097950475123
$ws.Range("D21").Value = 1.060502267677384
